$wb = $excel.ActiveWorkbook

# The existing "Sheet1" keeps its original data. We insert a brand-new
# worksheet named "New Values" in front of it, carrying the same first
# value plus the new CDR code.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$newSheet = $wb.Worksheets.Add($sheet1)
$newSheet.Name = "New Values"

$newSheet.Range("A1").Value = "value:1:1:1"
$newSheet.Range("A2").Value = "CDR.CLK0601322"

$newSheet.Activate()
